$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---

# Row 2: Dispositivo Raiz entry - IP and Interface changed
$ws.Cells.Item(2, 2).Value = "192.168.166.170"
$ws.Cells.Item(2, 3).Value = "FastEthernet2/0"

# Row 3: Vizinho Conectado entry - IP, Interface, Dispositivo, Plataforma changed
$ws.Cells.Item(3, 2).Value = "192.168.166.172"
$ws.Cells.Item(3, 3).Value = "FastEthernet1/0"
$ws.Cells.Item(3, 4).Value = "S2.lab.local"
$ws.Cells.Item(3, 5).Value = "Cisco 3725"

# Row 5: Dispositivo Raiz entry - IP, Dispositivo, Plataforma changed
$ws.Cells.Item(5, 2).Value = "192.168.166.173"
$ws.Cells.Item(5, 4).Value = "S3"
$ws.Cells.Item(5, 5).Value = "S3.lab.local"

# Row 6: Vizinho Conectado entry - IP, Interface, Dispositivo, Plataforma changed
$ws.Cells.Item(6, 2).Value = "192.168.166.172"
$ws.Cells.Item(6, 3).Value = "FastEthernet1/1"
$ws.Cells.Item(6, 4).Value = "S2.lab.local"
$ws.Cells.Item(6, 5).Value = "Cisco 3725"

# --- Add new rows 8 and 9 (row 7 stays blank as separator) ---

# Row 8: Dispositivo Raiz
$ws.Cells.Item(8, 1).Value = "Dispositivo Raiz"
$ws.Cells.Item(8, 2).Value = "192.168.166.172"
$ws.Cells.Item(8, 3).Value = "FastEthernet1/2"
$ws.Cells.Item(8, 4).Value = "S2"
$ws.Cells.Item(8, 5).Value = "S2.lab.local"

# Row 9: Vizinho Conectado
$ws.Cells.Item(9, 1).Value = "Vizinho Conectado"
$ws.Cells.Item(9, 2).Value = "192.168.166.174"
$ws.Cells.Item(9, 3).Value = "FastEthernet1/0"
$ws.Cells.Item(9, 4).Value = "S4.lab.local"
$ws.Cells.Item(9, 5).Value = "Cisco 3725"

# Row 10: blank separator row (mirrors rows 4 and 7). Rows 4/7 already
# contain empty (but present) cells, so copying one of them down materializes
# equally empty cells in row 10 instead of leaving the row absent entirely.
$ws.Range("A7:E7").Copy($ws.Range("A10:E10"))

$wb.Save()
